$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for I and J columns, copying H1's formatting (bold/border/center)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2-25: I = 1, J = value of H
for ($r = 2; $r -le 25; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
